# Apply the "changed % naive" update to the measles cost/benefit workbook.
# - Updates the Naive (C), Attack (D) and Vacc (E) columns for each DHB row
#   (rows 19-38) to reflect the latest MoH vaccine coverage rates.
# - The Population column (B) keeps its existing value but, like C/D/E,
#   drops its explicit "integer" number-format style so the cells fall back
#   to the default "Normal" style (matches the target workbook exactly).
# - Everything downstream (F..W, row 40 totals) is driven by formulas, so it
#   recalculates automatically once the inputs change.
# - Restores the view/selection state recorded for the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (C_new, D_new, E_new)
$rowData = @{
    19 = @(74410, 10662, 40320)
    20 = @(35150, 5041, 19056)
    21 = @(82195, 11770, 44525)
    22 = @(48368, 6928, 26204)
    23 = @(80484, 11633, 43820)
    24 = @(25934, 3730, 14082)
    25 = @(23864, 3478, 13053)
    26 = @(16854, 2439, 9182)
    27 = @(27397, 3855, 14697)
    28 = @(23502, 3398, 12799)
    29 = @(25902, 3719, 14051)
    30 = @(9444, 1344, 5099)
    31 = @(50047, 7025, 26811)
    32 = @(7473, 1077, 4063)
    33 = @(18842, 2728, 10268)
    34 = @(61138, 8730, 33067)
    35 = @(7005, 1002, 3793)
    36 = @(90522, 13170, 49463)
    37 = @(5554, 812, 3042)
    38 = @(10218, 1457, 5521)
}

foreach ($r in 19..38) {
    $vals = $rowData[$r]

    # Columns B:E lose their "s=2" (integer) style and become plain/"Normal".
    $ws.Range("B$r`:E$r").Style = "Normal"

    $ws.Range("C$r").Value = $vals[0]
    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
}

# Restore the recorded sheet selection/scroll position.
[void]$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B14").Select()

$excel.Calculate()
